$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Year of Release" column (C) stores its values as plain text (not
# numbers) in this sheet. Assigning a numeric-looking string via COM
# auto-converts it to a real number and also leaves a "quote prefix" style
# behind when forced to text with a leading apostrophe, so reset the style
# to Normal afterwards to keep the cell format identical to the source.
function Set-TextValue($rng, $value) {
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

# Row 75 / 76: swap Rome(2005) and Kota Factory(2019)
$ws.Range("B75").Value = "Kota Factory"
Set-TextValue $ws.Range("C75") "2019"
$ws.Range("B76").Value = "Rome"
Set-TextValue $ws.Range("C76") "2005"

# Row 119 / 120: swap Daredevil(2015) and Demon Slayer: Kimetsu no Yaiba(2019)
$ws.Range("B119").Value = "Demon Slayer: Kimetsu no Yaiba"
Set-TextValue $ws.Range("C119") "2019"
$ws.Range("B120").Value = "Daredevil"
Set-TextValue $ws.Range("C120") "2015"

# Row 122 / 123: swap Young Justice(2010) and Haikyuu!!(2014)
$ws.Range("B122").Value = "Haikyuu!!"
Set-TextValue $ws.Range("C122") "2014"
$ws.Range("B123").Value = "Young Justice"
Set-TextValue $ws.Range("C123") "2010"

# Row 163 / 164: swap Crash Landing on You(2019) and Jujutsu Kaisen(2020)
$ws.Range("B163").Value = "Jujutsu Kaisen"
Set-TextValue $ws.Range("C163") "2020"
$ws.Range("B164").Value = "Crash Landing on You"
Set-TextValue $ws.Range("C164") "2019"

# Row 171 / 172: swap Coupling(2000) and Mahabharat(1988)
$ws.Range("B171").Value = "Mahabharat"
Set-TextValue $ws.Range("C171") "1988"
$ws.Range("B172").Value = "Coupling"
Set-TextValue $ws.Range("C172") "2000"

# Row 180 / 181: swap Louie(2010) and Gullak(2019)
$ws.Range("B180").Value = "Gullak"
Set-TextValue $ws.Range("C180") "2019"
$ws.Range("B181").Value = "Louie"
Set-TextValue $ws.Range("C181") "2010"

# Row 208 / 209: swap Rurouni Kenshin: Trust and Betrayal(1999) and Detectorists(2014)
$ws.Range("B208").Value = "Detectorists"
Set-TextValue $ws.Range("C208") "2014"
$ws.Range("B209").Value = "Rurouni Kenshin: Trust and Betrayal"
Set-TextValue $ws.Range("C209") "1999"

# Row 227 / 228: swap John Adams(2008) and Erased(2016)
$ws.Range("B227").Value = "Erased"
Set-TextValue $ws.Range("C227") "2016"
$ws.Range("B228").Value = "John Adams"
Set-TextValue $ws.Range("C228") "2008"

# Row 242 / 243: swap Your Lie in April(2014) and I Love Lucy(1951)
$ws.Range("B242").Value = "I Love Lucy"
Set-TextValue $ws.Range("C242") "1951"
$ws.Range("B243").Value = "Your Lie in April"
Set-TextValue $ws.Range("C243") "2014"

# Row 251: replace Clannad: After Story(2008) with Avrupa Yakasi(2004)
$ws.Range("B251").Value = "Avrupa Yakasi"
Set-TextValue $ws.Range("C251") "2004"
